# Extend the "Julian" date-lookup tab (col A = MM.DD style date, col B =
# sequential day-of-year index) through the end of August and the first
# day of September, mirroring the tidy-script update that now also walks
# the GRG butterfly survey's later-season rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Julian")

# Column A values (dates) and column B values (running day count) that
# continue the existing sequence (row 88 was 8.14 / 226).
$dates = @(8.15, 8.16, 8.17, 8.18, 8.19, 8.20, 8.21, 8.22, 8.23, 8.24, 8.25, 8.26, 8.27, 8.28, 8.29, 8.30, 8.31, 9.01)
$startRow = 89
$startDay = 227

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $dateCell = $ws.Cells.Item($row, 1)
    $dayCell = $ws.Cells.Item($row, 2)

    $dateCell.Value = $dates[$i]
    $dayCell.Value = $startDay + $i

    # Match the formatting already used for column A / column B in this table
    # (copied from the existing row 88 style: numFmtId 2 "0.00", centered).
    $dateCell.NumberFormat = "0.00"
    $dateCell.HorizontalAlignment = -4108
    $dateCell.VerticalAlignment = -4108
    $dayCell.HorizontalAlignment = -4108
    $dayCell.VerticalAlignment = -4108
}

# Make "Julian" the active sheet/tab (was "Visual History") and leave the
# selection on the new last cell, the way Excel would after typing this in.
$ws.Activate()
$ws.Range("B106").Select() | Out-Null
